# Fruta / hortaliza, semanal
# The data rows (2-19) get their per-record values (Fecha, Volumen,
# Precio minimo/maximo/promedio, Origen, Precio $/Kg) reshuffled to
# correspond to a different weekly ordering. Columns D, J, K, L, M, O, P
# are affected; all other columns stay as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: target row -> source row (values currently sitting in source
# row get copied, unmodified, into target row).
$map = @{
    2  = 18
    3  = 8
    4  = 5
    5  = 19
    6  = 11
    7  = 12
    8  = 10
    9  = 7
    10 = 14
    11 = 6
    12 = 4
    13 = 16
    14 = 17
    15 = 2
    16 = 15
    17 = 3
    18 = 9
    19 = 13
}

# Columns whose values move together as a record.
$cols = @(4, 10, 11, 12, 13, 15, 16)   # D, J, K, L, M, O, P

# Snapshot the original values for every affected column/row before
# writing anything, since several rows are both sources and targets.
$snapshot = @{}
for ($r = 2; $r -le 19; $r++) {
    foreach ($c in $cols) {
        $snapshot["$r-$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

foreach ($targetRow in $map.Keys) {
    $sourceRow = $map[$targetRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($targetRow, $c).Value2 = $snapshot["$sourceRow-$c"]
    }
}
